$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of purely-numeric-looking data that must be stored as TEXT
# (matching the existing Account_ID / Customer_ID / PD columns), not as
# numbers. Writing numeric-looking strings straight into Range.Value
# makes Excel parse them as numbers (or, if forced via NumberFormat, tags
# them with a quote-prefix style) - neither matches the source data,
# which is plain shared-string text with no cell style override.
#
# Trick: build each value as a text FORMULA result in a scratch cell far
# away from the used range, then copy/paste-special *values only* into
# the destination. A paste of an already-text value keeps it text
# without forcing a quote-prefix cell style.

$data = @(
    @("118518", "1008784402", "17898937", "6020"),
    @("118518", "1008784413", "17898952", "6020"),
    @("118518", "1008784417", "17898955", "6020")
)

$startRow = 8
$helper = $ws.Range("Z1")

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    for ($c = 1; $c -le 4; $c++) {
        $text = $data[$i][$c - 1]
        $helper.Formula = "=""" + $text + """"
        $helper.Copy()
        $dst = $ws.Cells.Item($row, $c)
        $dst.PasteSpecial(-4163)  # xlPasteValues
    }
}

$helper.ClearContents()
$excel.CutCopyMode = 0
